# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets,
# continuing the existing "Date, Timestamp, Hour, Location, Value, Status"
# table layout. Date and percentage text values are forced to store as
# literal text (matching the existing text-log format) rather than being
# auto-converted by Excel into real dates/numbers.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($ws, $rowNum, $date, $timestamp, $hour, $location, $value, $status, $valueIsPercent)

    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $date

    $ws.Cells.Item($rowNum, 2).Value = $timestamp
    $ws.Cells.Item($rowNum, 3).Value = $hour
    $ws.Cells.Item($rowNum, 4).Value = $location

    if ($valueIsPercent) {
        $ws.Cells.Item($rowNum, 5).NumberFormat = "@"
    }
    $ws.Cells.Item($rowNum, 5).Value = $value

    $ws.Cells.Item($rowNum, 6).Value = $status
}

$logDate = "2026-02-06"

# --- PIR sheet: rows 326-339 (motion detection log) ---
$ws_pir = $wb.Worksheets.Item("PIR")
$pirRows = @(
    ("326","10:07:04","10:00","Bathroom","No Motion","Inactive"),
    ("327","10:07:05","10:00","Bathroom","Motion Detected","Active"),
    ("328","10:07:13","10:00","Bathroom","No Motion","Inactive"),
    ("329","10:07:13","10:00","Bathroom","Motion Detected","Active"),
    ("330","10:07:21","10:00","Bathroom","No Motion","Inactive"),
    ("331","10:07:23","10:00","Bathroom","Motion Detected","Active"),
    ("332","10:07:30","10:00","Bathroom","No Motion","Inactive"),
    ("333","10:07:31","10:00","Bathroom","Motion Detected","Active"),
    ("334","10:07:38","10:00","Bathroom","No Motion","Inactive"),
    ("335","10:07:39","10:00","Bathroom","Motion Detected","Active"),
    ("336","10:07:48","10:00","Bathroom","No Motion","Inactive"),
    ("337","10:07:49","10:00","Bathroom","Motion Detected","Active"),
    ("338","10:07:56","10:00","Bathroom","No Motion","Inactive"),
    ("339","10:07:58","10:00","Bathroom","Motion Detected","Active")
)
foreach ($row in $pirRows) {
    Add-LogRow $ws_pir $row[0] $logDate $row[1] $row[2] $row[3] $row[4] $row[5] $false
}

# --- Humidity sheet: rows 217-227 (humidity % log) ---
$ws_hum = $wb.Worksheets.Item("Humidity")
$humRows = @(
    ("217","10:07:11","10:00","Bathroom","69.5%","Active"),
    ("218","10:07:16","10:00","Bathroom","68.6%","Active"),
    ("219","10:07:21","10:00","Bathroom","69.5%","Active"),
    ("220","10:07:26","10:00","Bathroom","68.6%","Active"),
    ("221","10:07:31","10:00","Bathroom","69.7%","Active"),
    ("222","10:07:36","10:00","Bathroom","68.5%","Active"),
    ("223","10:07:41","10:00","Bathroom","69.4%","Active"),
    ("224","10:07:46","10:00","Bathroom","68.2%","Active"),
    ("225","10:07:51","10:00","Bathroom","69.2%","Active"),
    ("226","10:07:57","10:00","Bathroom","68.1%","Active"),
    ("227","10:08:01","10:00","Bathroom","69.2%","Active")
)
foreach ($row in $humRows) {
    Add-LogRow $ws_hum $row[0] $logDate $row[1] $row[2] $row[3] $row[4] $row[5] $true
}

# --- Temperature sheet: rows 216-227 (temperature log) ---
$ws_temp = $wb.Worksheets.Item("Temperature")
$tempRows = @(
    ("216","10:07:04","10:00","Bathroom","27.9C","Active"),
    ("217","10:07:12","10:00","Bathroom","27.9C","Active"),
    ("218","10:07:17","10:00","Bathroom","27.9C","Active"),
    ("219","10:07:22","10:00","Bathroom","27.9C","Active"),
    ("220","10:07:27","10:00","Bathroom","27.9C","Active"),
    ("221","10:07:32","10:00","Bathroom","28.0C","Active"),
    ("222","10:07:37","10:00","Bathroom","27.9C","Active"),
    ("223","10:07:42","10:00","Bathroom","28.0C","Active"),
    ("224","10:07:47","10:00","Bathroom","27.9C","Active"),
    ("225","10:07:52","10:00","Bathroom","28.0C","Active"),
    ("226","10:07:58","10:00","Bathroom","27.9C","Active"),
    ("227","10:08:02","10:00","Bathroom","28.0C","Active")
)
foreach ($row in $tempRows) {
    Add-LogRow $ws_temp $row[0] $logDate $row[1] $row[2] $row[3] $row[4] $row[5] $false
}
